$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; remove protection so the cells can be updated,
# then restore protection afterwards.
$ws.Unprotect()

# Update the confidentiality / "as of" date note (A13)
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-09 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for rows 2-10
$ws.Range("D2").Value = 0.1300689505693754
$ws.Range("E2").Value = 0.01238537572942722

$ws.Range("D3").Value = 0.1096687114482765
$ws.Range("E3").Value = -0.008212475974139544

$ws.Range("D4").Value = 0.110144907058881
$ws.Range("E4").Value = -0.001921008145074565

$ws.Range("D5").Value = 0.1193318733493163
$ws.Range("E5").Value = -0.004113629307668343

$ws.Range("D6").Value = 0.1196926130825006
$ws.Range("E6").Value = 0.002061288992716959

$ws.Range("D7").Value = 0.1505032036626603
$ws.Range("E7").Value = -0.005815563555801706

$ws.Range("D8").Value = 0.1302529134611984
$ws.Range("E8").Value = -0.002548853016142827

$ws.Range("D9").Value = 0.1303368273677915
$ws.Range("E9").Value = -0.006332425403567044

$ws.Range("E10").Value = -0.001778058837873742

# Restore sheet protection to match the original protected state.
$ws.Protect()
